$d = $word.ActiveDocument

# Helper: replace the first occurrence of $OldText found at/after character
# position $StartPos with $NewText. Returns the end position of the
# replacement (so callers can keep advancing through the document in
# order without re-matching earlier, already-handled occurrences).
function Replace-NextOccurrence($OldText, $NewText, $StartPos) {
    $docEnd = $d.Content.End
    $rng = $d.Range($StartPos, $docEnd)
    $found = $rng.Find.Execute($OldText, $true, $false, $false, $false, $false, `
                                $true, 0, $false, $NewText, 1)
    if (-not $found) {
        throw "Replace-NextOccurrence: could not find '$OldText' at/after $StartPos"
    }
    return $rng.End
}

$cursor = 0

# Organizzazione: ROSSI AMAZZONIA -> GIUBILEO SPA
$cursor = Replace-NextOccurrence "ROSSI AMAZZONIA" "GIUBILEO SPA" $cursor

# Data/e di Verifica: 23-24-25-30-31 ottobre 2023 e 2-3-6 -> 23-25 DICEMBRE 2025
$cursor = Replace-NextOccurrence "23-24-25-30-31 ottobre 2023 e 2-3-6" "23-25 DICEMBRE 2025" $cursor

# SITO 2: TERNI -> SITO 2: GENZANO   (only this exact run; the later
# "SITO2 -  TERNI - ON SITE --" run must stay untouched)
$cursor = Replace-NextOccurrence "TERNI " "GENZANO" $cursor

# Auditor name: MARIANGELO GIOVANNINI -> Marianna Pezzuca
$cursor = Replace-NextOccurrence "MARIANGELO GIOVANNINI" "Marianna Pezzuca" $cursor

# Five "Attività del:" dates, in document order, each distinct:
$cursor = Replace-NextOccurrence "24/1/2025" "23.10.2023" $cursor
$cursor = Replace-NextOccurrence "24/1/2025" "23.10.2023 " $cursor
$cursor = Replace-NextOccurrence "24/1/2026" "24.10.2023" $cursor
$cursor = Replace-NextOccurrence "24/1/2027" " 25.10.202" $cursor
$cursor = Replace-NextOccurrence "24/1/2028" "30.10.2023" $cursor

Write-Host "Done. Final cursor:" $cursor
